# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 205) from 2023-09-13 (serial 45182) to 2023-09-15
# (serial 45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 205; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
